# Remove the "Gewinn 2023" entry from the Einnahmen (income) sheet.
# Commit message: "Gewinn aus dem letzten Jahr aus den Einnahmen entfernt"
# (Removed last year's profit from the income sheet)

$wb = $excel.ActiveWorkbook

# Work on the "Einnahmen" sheet where the "Gewinn 2023" row lives (row 8:
# Kategorie=Sonstiges, Bezeichnung=Gewinn 2023, Datum=45322, Betrag=4001.43,
# Firmennamen=Atelierkino).
$wsEin = $wb.Worksheets.Item("Einnahmen")
$wsEin.Activate()

# Select the whole row (mirrors the user selecting row 8 before deleting it)
# and delete it, which shifts all following rows up by one and shrinks the
# table / autofilter / sort state / data validation ranges accordingly.
$row = $wsEin.Rows.Item(8)
$row.Select()
$row.Delete()

$wb.Save()
